# Update the PROSOUL work-estimate sheet: flip the in-progress rows'
# STATUS from PENDING to PROGRESS, and leave the selection where the
# author left it (C4:C5) after making the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$progressCells = @("H3", "H4", "H5", "H13", "H14", "H15", "H16", "H17", "H18", "H19")
foreach ($cellRef in $progressCells) {
    $ws.Range($cellRef).Value = "PROGRESS"
}

$ws.Range("C4:C5").Select()

$wb.Save()
